$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45180 -> 45181) for every data row, from row 2 through row 533.
$lastRow = 533
$ws.Range("C2:C$lastRow").Value = 45181
